# Scene.xlsx reconstitution: remove the unused "PioneerNoob" path string and
# populate the (previously blank) FilePath column (B) for every scene row
# with "../../NFDataCfg/Ini/Scene/<ID>.xml", matching the sibling *.xml ini
# files. Also moves the saved selection/view back to the default top-left
# with the cursor parked on B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ID 1) already had a FilePath value (the stale PioneerNoob path) -
# just replace its text; its existing style/format is left untouched.
$ws.Range("B2").Value = "../../NFDataCfg/Ini/Scene/1.xml"

# Rows 3-7 (IDs 2-6) had an empty FilePath cell. Fill them in and give them
# the distinct font (same 11pt black Song typeface, family 3) + text number
# format ("@") that the new values use.
$rows = @(3, 4, 5, 6, 7)
$ids  = @(2, 3, 4, 5, 6)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $id = $ids[$i]
    $cell = $ws.Range("B$r")
    $cell.Value = "../../NFDataCfg/Ini/Scene/$id.xml"
    $cell.NumberFormat = "@"
    $cell.Font.Name = "宋体"
    $cell.Font.Size = 11
    $cell.Font.Color = 0
    $cell.Font.Family = 3
}

# Restore the view to the default top-left cell and park the selection on B5.
[void]$ws.Range("B5").Select()
